$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "SKU"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Manufacturer"
$ws.Range("D1").Value = "ManufacturerPart"
$ws.Range("E1").Value = "ProcessRequest"
$ws.Range("F1").Value = "SortingRequest"
$ws.Range("G1").Value = "Unit"
$ws.Range("H1").Value = "UnitPrice"
$ws.Range("I1").Value = "Currency"
$ws.Range("J1").Value = "Qty"
